$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Cash"
$ws.Range("A3").Value = "US Equities"
$ws.Range("A4").Value = "European Equities"
$ws.Range("A5").Value = "EU High Yield"
$ws.Range("A6").Value = "EU Corporate"
$ws.Range("A7").Value = "Greek Gov"
$ws.Range("A8").Value = "Euro Gov"
